# Update the NATMI LR-pairs sheet (Bmp7-Bmpr1b) with refreshed TPM numbers.
# Rows 2-5 get new sending/target-cluster labels and recalculated metrics,
# and two new rows (6-7, MuSCs sending cluster) are appended.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bmp7"
$ws.Range("C2").Value = "Bmpr1b"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03443933333333333
$ws.Range("H2").Value = 0.103318
$ws.Range("I2").Value = 0.05823261822459219
$ws.Range("J2").Value = 0.0582326182245922
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.707786666666667
$ws.Range("N2").Value = 5.12336
$ws.Range("O2").Value = 0.7091726973716084
$ws.Range("P2").Value = 0.7091726973716084
$ws.Range("Q2").Value = 0.05881503427555555
$ws.Range("R2").Value = 0.52933530848
$ws.Range("S2").Value = 0.04129698294134512
$ws.Range("T2").Value = 0.04129698294134513

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bmp7"
$ws.Range("C3").Value = "Bmpr1b"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03443933333333333
$ws.Range("H3").Value = 0.103318
$ws.Range("I3").Value = 0.05823261822459219
$ws.Range("J3").Value = 0.0582326182245922
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7003526666666667
$ws.Range("N3").Value = 2.101058
$ws.Range("O3").Value = 0.2908273026283917
$ws.Range("P3").Value = 0.2908273026283917
$ws.Range("Q3").Value = 0.02411967893822222
$ws.Range("R3").Value = 0.217077110444
$ws.Range("S3").Value = 0.01693563528324707
$ws.Range("T3").Value = 0.01693563528324707

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Bmp7"
$ws.Range("C4").Value = "Bmpr1b"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.51625
$ws.Range("H4").Value = 1.54875
$ws.Range("I4").Value = 0.8729143757654733
$ws.Range("J4").Value = 0.8729143757654734
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.707786666666667
$ws.Range("N4").Value = 5.12336
$ws.Range("O4").Value = 0.7091726973716084
$ws.Range("P4").Value = 0.7091726973716084
$ws.Range("Q4").Value = 0.8816448666666666
$ws.Range("R4").Value = 7.9348038
$ws.Range("S4").Value = 0.6190470424360544
$ws.Range("T4").Value = 0.6190470424360545

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Bmp7"
$ws.Range("C5").Value = "Bmpr1b"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.51625
$ws.Range("H5").Value = 1.54875
$ws.Range("I5").Value = 0.8729143757654733
$ws.Range("J5").Value = 0.8729143757654734
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7003526666666667
$ws.Range("N5").Value = 2.101058
$ws.Range("O5").Value = 0.2908273026283917
$ws.Range("P5").Value = 0.2908273026283917
$ws.Range("Q5").Value = 0.3615570641666667
$ws.Range("R5").Value = 3.2540135775
$ws.Range("S5").Value = 0.2538673333294189
$ws.Range("T5").Value = 0.253867333329419

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Bmp7"
$ws.Range("C6").Value = "Bmpr1b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.04072033333333334
$ws.Range("H6").Value = 0.122161
$ws.Range("I6").Value = 0.06885300600993445
$ws.Range("J6").Value = 0.06885300600993445
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.707786666666667
$ws.Range("N6").Value = 5.12336
$ws.Range("O6").Value = 0.7091726973716084
$ws.Range("P6").Value = 0.7091726973716084
$ws.Range("Q6").Value = 0.06954164232888889
$ws.Range("R6").Value = 0.62587478096
$ws.Range("S6").Value = 0.04882867199420877
$ws.Range("T6").Value = 0.04882867199420877

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Bmp7"
$ws.Range("C7").Value = "Bmpr1b"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.04072033333333334
$ws.Range("H7").Value = 0.122161
$ws.Range("I7").Value = 0.06885300600993445
$ws.Range("J7").Value = 0.06885300600993445
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7003526666666667
$ws.Range("N7").Value = 2.101058
$ws.Range("O7").Value = 0.2908273026283917
$ws.Range("P7").Value = 0.2908273026283917
$ws.Range("Q7").Value = 0.02851859403755556
$ws.Range("R7").Value = 0.256667346338
$ws.Range("S7").Value = 0.02002433401572568
$ws.Range("T7").Value = 0.02002433401572568
